$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B6").Value = 99
$ws.Range("B7").Value = 25

$ws.Range("A8").Value = "Total distance covered (km)"
$ws.Range("A9").Value = "Total energy consumption(WH/KM)"
$ws.Range("A10").Value = "Total SOC consumed(%)"

$ws.Range("A12").Value = "Peak Power(kW)"
$ws.Range("A13").Value = "Average Power(kW)"
$ws.Range("A14").Value = "Total Energy Regenerated(kWh)"

$ws.Range("A15").Value = "Regenerative Effectiveness(%)"
$ws.Range("B15").Value = 0.0008349768799234199

$ws.Range("A16").Value = "Highest Cell Voltage(V)"
$ws.Range("B16").Value = 3.357

$ws.Range("A17").Value = "Lowest Cell Voltage(V)"
$ws.Range("B17").Value = 2.988

$ws.Range("A18").Value = "Difference in Cell Voltage(V)"

$ws.Range("A19").Value = "Minimum Temperature(C)"
$ws.Range("A20").Value = "Maximum Temperature(C)"

$ws.Range("A21").Value = "Difference in Temperature(C)"
$ws.Range("B21").Value = 15

$ws.Range("A22").Value = "Maximum Fet Temperature-BMS(C)"
$ws.Range("A23").Value = "Maximum Afe Temperature-BMS(C)"
$ws.Range("A24").Value = "Maximum PCB Temperature-BMS(C)"
$ws.Range("A25").Value = "Maximum MCU Temperature(C)"
$ws.Range("A26").Value = "Maximum Motor Temperature(C)"
$ws.Range("A27").Value = "Abnormal Motor Temperature Detected(C)"

$ws.Range("A28").Value = "highest cell temp(C)"
$ws.Range("A29").Value = "lowest cell temp(C)"

$ws.Range("A30").Value = "Difference between Highest and Lowest Cell Temperature at 100% SOC(C)"

$ws.Range("A31").Value = "Battery Voltage(V)"
$ws.Range("B31").Value = 53

$ws.Range("A32").Value = "Total energy charged(kWh)"
$ws.Range("B32").Value = 1.468598774166667

$ws.Range("A33").Value = "Electricity consumption units(kW)"
$ws.Range("B33").Value = 0.0000001434402615805855

$ws.Range("A34").Value = "Idling time percentage"
$ws.Range("B34").Value = 37.44252152830031

$ws.Range("A35").Value = "Time spent in 0-10 km/h"
$ws.Range("B35").Value = 3.072485578128919

$ws.Range("A36").Value = "Time spent in 10-20 km/h"
$ws.Range("B36").Value = 3.206253657720926

$ws.Range("A37").Value = "Time spent in 20-30 km/h"
$ws.Range("B37").Value = 5.09154752947078

$ws.Range("A38").Value = "Time spent in 30-40 km/h"
$ws.Range("B38").Value = 6.993562411169635

$ws.Range("A39").Value = "Time spent in 40-50 km/h"
$ws.Range("B39").Value = 8.067887300392943

$ws.Range("A40").Value = "Time spent in 50-60 km/h"
$ws.Range("B40").Value = 7.779449878772677

$ws.Range("A41").Value = "Time spent in 60-70 km/h"
$ws.Range("B41").Value = 12.79157261098571

$ws.Range("A42").Value = "Time spent in 70-80 km/h"
$ws.Range("B42").Value = 15.4293119304406

$ws.Range("A43").Value = "Time spent in 80-90 km/h"
$ws.Range("B43").Value = 0
